# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column G holds the "K" (strikeouts) values for rows 2-31.
# These are recalculated values replacing the previous Strike# derived figures.
$kValues = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 4
    6  = 3
    7  = 3
    8  = 4
    9  = 3
    10 = 3
    11 = 3
    12 = 3
    13 = 1
    14 = 1
    15 = 9
    16 = 0
    17 = 1
    18 = 3
    19 = 3
    20 = 3
    21 = 3
    22 = 0
    23 = 3
    24 = 1
    25 = 4
    26 = 4
    27 = 3
    28 = 5
    29 = 1
    31 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
